$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.364.38"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "2.284.31"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'503.41"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("D6").Value = "'129.71"
$ws.Range("E6").Value = "  +1.62%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("E9").Value = "  +1.93%  "

$ws.Range("E10").Value = "  +1.33%  "

$ws.Range("E11").Value = "  +3.39%  "

$ws.Range("E12").Value = "  +1.70%  "

$ws.Range("D13").Value = "2.693.41"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").Value = "'23.07"
$ws.Range("E14").Value = "  +6.66%  "

$ws.Range("D15").Value = "54.334.53"
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").Value = "2.303.39"
$ws.Range("E17").Value = "  +0.87%  "

$ws.Range("D18").Value = "'10.28"
$ws.Range("E18").Value = "  +3.62%  "

$ws.Range("E19").Value = "  +2.40%  "

$ws.Range("D20").Value = "'305.55"
$ws.Range("E20").Value = "  +2.28%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").Value = "'62.11"
$ws.Range("E23").Value = "  -2.65%  "

$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("E25").Value = "  +2.29%  "

$ws.Range("E26").Value = "  +3.14%  "

$ws.Range("D27").Value = "'174.45"
$ws.Range("E27").Value = "  +6.71%  "

$ws.Range("E28").Value = "  +1.26%  "

$ws.Range("D29").Value = "'6.02"
$ws.Range("E29").Value = "  +2.79%  "

$ws.Range("E30").Value = "  +1.90%  "

$ws.Range("D31").Value = "'1.08"
$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "'17.82"
$ws.Range("E33").Value = "  +1.95%  "

$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").Value = "'0.944"
$ws.Range("E35").Value = "  +8.46%  "

$ws.Range("E36").Value = "  +1.85%  "

$ws.Range("E37").Value = "  +4.06%  "

$ws.Range("D38").Value = "'0.375"
$ws.Range("E38").Value = "  -0.43%  "

$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("E40").Value = "  +2.03%  "

$ws.Range("D41").Value = "'4.82"
$ws.Range("E41").Value = "  -0.55%  "

$ws.Range("D42").Value = "'125.18"
$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("D43").Value = "'0.0496"
$ws.Range("E43").Value = "  +3.35%  "

$ws.Range("E45").Value = "  +0.53%  "

$ws.Range("D46").Value = "'241.17"
$ws.Range("E46").Value = "  +1.10%  "

$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("E48").Value = "  +1.67%  "

$ws.Range("E49").Value = "  +0.95%  "

$ws.Range("D50").Value = "'16.41"
$ws.Range("E50").Value = "  +0.51%  "

$ws.Range("E51").Value = "  +0.13%  "

